# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 86
$wsExpo.Range("F3").Value = 4030
$wsExpo.Range("F4").Value = 2360
$wsExpo.Range("F9").Value = 198
$wsExpo.Range("F11").Value = 71
$wsExpo.Range("F12").Value = 130
$wsExpo.Range("F13").Value = 1502
$wsExpo.Range("F14").Value = 266
$wsExpo.Range("F15").Value = 2843
$wsExpo.Range("F16").Value = 195

# --- Sheet "全部类型" (sheet4) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 86
$wsAll.Range("F3").Value = 4030
$wsAll.Range("F4").Value = 2360
$wsAll.Range("F10").Value = 198
$wsAll.Range("F12").Value = 71
$wsAll.Range("F13").Value = 130
$wsAll.Range("F16").Value = 1502
$wsAll.Range("F17").Value = 266
$wsAll.Range("F18").Value = 2843
$wsAll.Range("F19").Value = 195

$wb.Save()
